# Update the communication style dictionary (Sheet1)
# Adds talkative/informality/sentimentality/conciseness/conversational-dominance
# rows, tweaks the "authoritative" survey item text, and reformats the table
# (header font, wrapped/Calibri body font, row heights, column C width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10
$rsquo = [char]8217

# ---------------------------------------------------------------------------
# 1. Cell values
# ---------------------------------------------------------------------------

# Header row (unchanged text, gets re-styled below)
$ws.Range("A1").Value = "style"
$ws.Range("B1").Value = "definition"
$ws.Range("C1").Value = "survey_item"

# Row 2 - authoritative (survey_item text tweaked: "not very likely" -> "very likely")
$ws.Range("A2").Value = "authoritative"
$ws.Range("B2").Value = "Authoritative is the tendency to command or demand others in a conversation, without considering the others' willingness and concerns."
$ws.Range("C2").Value = "I am very likely to tell someone what they should do; I sometimes insist that otheres do what I say; I expect people to obey when I ask them to do something; When I feel others should do something for me, I ask for it in a demanding tone of voice. "

# Row 3 - talkative (new)
$ws.Range("A3").Value = "talkative"
$ws.Range("B3").Value = "Talkativeness is a tendency to initiate a conversation, talk a lot, and avoid silence in a conversation."
$ws.Range("C3").Value = "I always have a lot to say; I have a hard time keeping myself silent when around other people; I am always the one who breaks a silence by starting to talk; I like to talk a lot."

# Row 4 - informality (new)
$ws.Range("A4").Value = "informality"
$ws.Range("B4").Value = "Informality is a tendency to talk casually and avoid being formal, distant, or stiff in a conversation."
$ws.Range("C4").Value = "I never communicate with others in a distant manner; I never behave formally when I meet someone; I always address others in a very casual way; I never come across as somewhat stiff when dealing with people." + $nl

# Row 5 - sentimentality (new)
$ws.Range("A5").Value = "sentimentality"
$ws.Range("B5").Value = "Sentimentality is a tendency to express one's own emotions or display empathic emotional responses to others in a conversation."
$ws.Range("C5").Value = "When I see others cry, I have difficulty holding back my tears; During a conversation, I am easily overcome by emotions; When describing my memories, i sometimes get visibily emotional; People can tell that I am emotionally touched by some topics of conversation."

# Row 6 - conciseness (new)
$ws.Range("A6").Value = "conciseness"
$ws.Range("B6").Value = "Conciseness is the tendency to use as few words as possible to clearly convey ideas and explain things in a conversation, and avoid being long-winded."
$ws.Range("C6").Value = "I don" + $rsquo + "t need a lot of words to get my message across; Most of the time, I only need a few words to explain something; With a few words I can usually clarify my point to everybody."

# Row 7 - conversational dominance (new)
$ws.Range("A7").Value = "conversational dominance"
$ws.Range("B7").Value = "Conversational dominance is the tendency to take the lead in a conversation and detremine its topics and directions."
$ws.Range("C7").Value = "I often take the lead in a conversation; I often determine which topics are talked about during a conversation; I often determine the direction of a conversation."

# ---------------------------------------------------------------------------
# 2. Column width / row heights
# ---------------------------------------------------------------------------

$ws.Columns.Item(3).ColumnWidth = 83.33203125

$ws.Rows.Item(2).RowHeight = 169
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 68
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 34
$ws.Rows.Item(7).RowHeight = 34

# ---------------------------------------------------------------------------
# 3. Fonts / alignment
#    - header row + informality..conversational-dominance (A/B): Calibri 12
#    - authoritative/talkative (A/B): Calibri 12
#    - survey_item column: Calibri 12 + wrap text
#    Font.ThemeColor is reasserted everywhere (not just on the previously
#    red Menlo cells) so every touched cell ends up on the same "automatic
#    text colour" font rather than inheriting a stray override.
#    NOTE: comma-joined multi-area Range() strings only affect the first
#    area in this host, so each area is touched with its own Range() call.
# ---------------------------------------------------------------------------

foreach ($a in @("A1:C1", "A4:B4", "A5:B5", "A6:B6", "A7:B7")) {
    $r = $ws.Range($a)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 12
    $r.Font.ThemeColor = 1
}

foreach ($a in @("A2:B2", "A3:B3")) {
    $r = $ws.Range($a)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 12
    $r.Font.ThemeColor = 1
}

foreach ($a in @("C2", "C3", "C4", "C5", "C6", "C7")) {
    $r = $ws.Range($a)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 12
    $r.Font.ThemeColor = 1
    $r.WrapText = $true
}

# ---------------------------------------------------------------------------
# 4. Selection (cosmetic, matches author's final cursor position)
# ---------------------------------------------------------------------------

$ws.Range("C10").Select()
